$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (FechaRecepcion) and Column S (FechaDevolucion) in row 2
# change from a Date value to a Text value "1/4/2022"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "1/4/2022"

$ws.Cells.Item(2, 19).NumberFormat = "@"
$ws.Cells.Item(2, 19).Value = "1/4/2022"
